$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 text with new conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 13.16 = 53368.68 pesos`n✅ 53368.68 pesos = 13.08 = 971.89 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update tasas sheet values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 76
$wsTasas.Range("O10").Value = 4056.02
$wsTasas.Range("N12").Value = 4080
$wsTasas.Range("O12").Value = 74.3
